# Insert a new data row at row 204 (pushing all existing rows from 204
# downward down by one, which also grows the used range from R332 to R333),
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("204:204").Insert()

$ws.Range("A204").Value = 3
$ws.Range("B204").Value = "Femacal de La Calera"
$ws.Range("C204").Value = "Coquimbo"
$ws.Range("D204").Value = 44606
$ws.Range("E204").Value = 5
$ws.Range("F204").Value = 100112017
$ws.Range("G204").Value = "Apio"
$ws.Range("H204").Value = "Americana (o)"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 120
$ws.Range("K204").Value = 9000
$ws.Range("L204").Value = 9000
$ws.Range("M204").Value = 9000
$ws.Range("N204").Value = "$/docena de matas"
$ws.Range("O204").Value = "Región Metropolitana"
$ws.Range("P204").Value = 1500
$ws.Range("Q204").Value = 6
$ws.Range("R204").Value = "Hortaliza"
